$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 135, shifting existing rows 135-218 down to 136-219
$ws.Rows("135:135").Insert()

# Fill in the values for the newly inserted row 135 (a new weekly garlic price record)
$ws.Range("A135").Value = 8
$ws.Range("B135").Value = "Terminal La Palmera de La Serena"
$ws.Range("C135").Value = "Coquimbo"
$ws.Range("D135").Value = 44603
$ws.Range("E135").Value = 4
$ws.Range("F135").Value = 100112003
$ws.Range("G135").Value = "Ajo"
$ws.Range("H135").Value = "Chino"
$ws.Range("I135").Value = "Primera"
$ws.Range("J135").Value = 540
$ws.Range("K135").Value = 18500
$ws.Range("L135").Value = 19000
$ws.Range("M135").Value = 18750
$ws.Range("N135").Value = "$/caja 10 kilos"
$ws.Range("O135").Value = "China"
$ws.Range("P135").Value = 1875
$ws.Range("Q135").Value = 10
$ws.Range("R135").Value = "Hortaliza"
